$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Sheet1"): delete row 1 (the "up73h59u" entry), shifting all rows up ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows.Item(1).Delete()

# --- Sheet2 ("used"): append a new row with the "up73h59u" entry now marked as used ---
$ws2 = $wb.Worksheets.Item("used")
$ws2.Range("A20").Value = "up73h59u"
$ws2.Range("B20").Value = "ChatGPT Image 2026年1月18日 10_33_10.png"
$ws2.Range("C20").Value = "2026-01-18 10:35:13"
